$d = $word.ActiveDocument

$replacements = @(
    @("892×6=5352", "226×5=1130"),
    @("132×2=264", "659×7=4613"),
    @("729×3=2187", "187×8=1496"),
    @("126×7=882", "643×6=3858"),
    @("393×8=3144", "879×6=5274"),
    @("660×7=4620", "423×4=1692"),
    @("417×8=3336", "529×2=1058"),
    @("271×4=1084", "697×5=3485"),
    @("892×4=3568", "299×4=1196"),
    @("444×4=1776", "664×7=4648"),
    @("363×2=726", "336×5=1680"),
    @("326×4=1304", "924×4=3696"),
    @("185×5=925", "836×8=6688"),
    @("370×2=740", "101×6=606"),
    @("613×7=4291", "973×6=5838"),
    @("461×9=4149", "333×7=2331"),
    @("964×5=4820", "376×7=2632"),
    @("123×7=861", "589×2=1178"),
    @("493×5=2465", "486×6=2916"),
    @("284×7=1988", "888×7=6216"),
    @("389×6=2334", "925×4=3700"),
    @("146×3=438", "652×6=3912"),
    @("463×6=2778", "967×6=5802"),
    @("183×2=366", "365×5=1825"),
    @("517×7=3619", "991×2=1982")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
